# Automatische test-sync: 2025-08-03 18:51:50
# Append a new log row (#46) to the "Logs" sheet and update the
# "Dashboard" summary count for "Inkoop / Bestellingen".

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 46

$logs.Cells.Item($newRow, 1).Value = "Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #18: Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Cells.Item($newRow, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-03 18:51:06"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Bump the Dashboard summary count for "Inkoop / Bestellingen" (row 4, col B).
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(4, 2).Value = $dashboard.Cells.Item(4, 2).Value2 + 1

# Extend the conditional-formatting ranges (D/G/H/I/J) from row 45 to row 46
# so the newly appended row keeps getting highlighted like the rest.
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "45")
    $newRange = $logs.Range($col + "2:" + $col + "46")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
